# Complete order from home page into category details
# Append two new data rows (Nathan Wagner, Reamer Schickowski) below the
# existing rows on Sheet1, extending the used range from A1:D3 to A1:D5.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A4").Value = "Nathan"
$ws.Range("B4").Value = "Wagner"
$ws.Range("C4").Value = "Test@1234"
$ws.Range("D4").Value = "NathanWagner@yopmail.com"

$ws.Range("A5").Value = "Reamer"
$ws.Range("B5").Value = "Schickowski"
$ws.Range("C5").Value = "Test@1234"
$ws.Range("D5").Value = "ReamerSchickowski@yopmail.com"
